# Applies the "prueba_2050" milestone-year column (G) edits to TimePeriods,
# drops the stale 2019-era selection / active sheet, and bumps the discount
# rate on Constants per the commit message
# ("Se encontro el error de reservas, algunas estan por encima de 2050").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: TimePeriods  (milestone years table, new column G "prueba_2050")
# ---------------------------------------------------------------------
$tp = $wb.Worksheets.Item("TimePeriods")

# New, slightly wider column F (matches the sibling columns visually)
$tp.Columns.Item(6).ColumnWidth = 13

# New header cell G35, formatted like its left neighbour F35
$tp.Range("G35").Value = "prueba_2050"
$tp.Range("F35").Copy()
[void]$tp.Range("G35").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New "endyear" values in column G, mirroring column E but capped at 2050
$tp.Range("G37").Value = 2019
$tp.Range("G38").Value = 2020
$tp.Range("G39").Value = 2023
$tp.Range("G40").Value = 2025
$tp.Range("G41").Value = 2027
$tp.Range("G42").Value = 2030
$tp.Range("G43").Value = 2033
$tp.Range("G44").Value = 2035
$tp.Range("G45").Value = 2037
$tp.Range("G46").Value = 2040
$tp.Range("G47").Value = 2043
$tp.Range("G48").Value = 2045
$tp.Range("G49").Value = 2047
$tp.Range("G50").Value = 2050

$tp.Range("G51").Select()

# ---------------------------------------------------------------------
# Sheet: Interpol_Extrapol_Defaults (loses tab focus, selection untouched)
# ---------------------------------------------------------------------
$ied = $wb.Worksheets.Item("Interpol_Extrapol_Defaults")
$ied.Range("B9").Select()

# ---------------------------------------------------------------------
# Sheet: Constants (discount rate fix + becomes the active tab)
# ---------------------------------------------------------------------
$cst = $wb.Worksheets.Item("Constants")
$cst.Range("E8").Value = 0.12
$cst.Range("E9").Select()
